$d = $word.ActiveDocument

# 1) Update the letter date.
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "September 21, 2025", 2)

# 2) Split the mailing address line "2113 Wendover Ln, San Jose CA 95121" into
#    "2113 Wendover Ln" / "San Jose, CA 95121" on their own paragraphs, followed
#    by a new blank paragraph (the pre-existing blank paragraph that used to
#    follow the address line is left untouched after these new ones).
$d.Content.Find.Execute("2113 Wendover Ln, San Jose CA 95121", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2113 Wendover Ln" + [char]13 + "San Jose, CA 95121" + [char]13, 2)

# 3) Remove the two blank paragraphs that used to sit right after the
#    "Board of Directors" signature line (a blank "No Spacing" paragraph and a
#    blank "Title" paragraph), keeping the following blank "Title" paragraph.
$all = $d.Paragraphs
for ($i = 1; $i -le $all.Count; $i++) {
    $p = $all.Item($i)
    $ptext = $p.Range.Text.TrimEnd([char]13)
    if ($ptext -eq "Vietnam Town Condominium Owners Association Board of Directors") {
        $idx = $i
        break
    }
}

$target = $all.Item($idx + 1)
$target.Range.Delete()

$all2 = $d.Paragraphs
$target2 = $all2.Item($idx + 1)
$target2.Range.Delete()
